$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-style pattern: for "Price" (column D) cells whose new value would
# otherwise be auto-parsed by Excel as a real number (losing the original
# text formatting, e.g. trailing zeros or the site's dot-grouped display
# style), force the cell to Text first so the literal string is preserved.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.024.42"
$ws.Range("E2").Value = "  +0.87%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.643.12"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.38%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.06"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.81%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.35%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.57%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0639"
$ws.Range("E9").Value = "  +1.28%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  +0.67%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.89%  "

# Row 12 - now WrappedEther (was WrappedliquidstakedEther2.0)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.693.49"
$ws.Range("E12").Value = "  +4.39%  "

# Row 13 - now WrappedliquidstakedEther2.0 (was WrappedEther)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.871.08"
$ws.Range("E13").Value = "  +0.95%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.88%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.06%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  +1.16%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.46"
$ws.Range("E17").Value = "  +1.59%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "26.107.88"
$ws.Range("E18").Value = "  +1.11%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.41%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.54"
$ws.Range("E20").Value = "  +1.61%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.12%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  -0.08%  "

# Row 24 - Stellar
$ws.Range("E24").Value = "  +4.50%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +0.15%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.38%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "143.29"
$ws.Range("E27").Value = "  +0.30%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.89"
$ws.Range("E28").Value = "  +1.06%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.92%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.36%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0497"
$ws.Range("E31").Value = "  +0.15%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.73%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.60%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -2.30%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.92%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +0.89%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.129.84"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38 - ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.540"
$ws.Range("E38").Value = "  -1.11%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -0.26%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.79%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.48"
$ws.Range("E41").Value = "  +0.82%  "

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.32"
$ws.Range("E42").Value = "  +0.30%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -0.45%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.780.63"
$ws.Range("E44").Value = "  +0.98%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0118"
$ws.Range("E45").Value = "  +4.70%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.75"
$ws.Range("E46").Value = "  +1.33%  "

# Row 47 - Cronos
$ws.Range("E47").Value = "  -0.32%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +1.02%  "

# Row 49 - EnergySwap
$ws.Range("E49").Value = "  +2.59%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.18%  "

# Row 51 - Algorand
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -0.29%  "
